$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-unused rows (old rows 6..12 covering motherboard..SSD)
$ws.Rows("6:12").Delete()

# Correct the migrated values for the remaining rows
$ws.Range("C2").Value = 100000
$ws.Range("D2").Value = 1

$ws.Range("C3").Value = 12000
$ws.Range("D3").Value = 5
# E3 was "false" (wrong migration); fix it to the "true" text value already
# used elsewhere in the sheet (copy as value so it reuses the shared string
# instead of becoming a native boolean).
$ws.Range("E2").Copy()
$ws.Range("E3").PasteSpecial(-4163)

$ws.Range("C4").Value = 15000

$ws.Range("C5").Value = 3200
$ws.Range("D5").Value = 12

# Header row no longer carries the centered style
$ws.Range("A1:E1").ClearFormats()

# Data cells switch from centered to right-aligned
$ws.Range("A2:E5").HorizontalAlignment = -4152

# Update the view: zoom + selection
$excel.ActiveWindow.Zoom = 175
$ws.Range("D7").Select() | Out-Null
